# RHODE_ISLAND_2022.xlsx cleanup:
#  - rename header columns to snake_case field names
#  - normalize capitalization of Spanish connector words ("de", "del", "el", "la")
#    to title case ("De", "Del", "El", "La") inside municipality/state names
#  - drop the trailing footnote/metadata rows (114-118), shrinking the used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalization fixes on individual data cells ---
$ws.Range("A6").Value  = "Ciudad De México"
$ws.Range("A14").Value = "Estado De México"
$ws.Range("B14").Value = "Ecatepec De Morelos"
$ws.Range("B16").Value = "Naucalpan De Juárez"
$ws.Range("B19").Value = "Tenango Del Valle"
$ws.Range("B20").Value = "Tlalnepantla De Baz"
$ws.Range("B22").Value = "Apaseo El Alto"
$ws.Range("B23").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B28").Value = "Acapulco De Juárez"
$ws.Range("B31").Value = "Ayutla De Los Libres"
$ws.Range("B33").Value = "Chilpancingo De Los Bravo"
$ws.Range("B34").Value = "Coyuca De Benítez"
$ws.Range("B36").Value = "Zihuatanejo De Azueta"
$ws.Range("B47").Value = "Huasca De Ocampo"
$ws.Range("B49").Value = "Mixquiahuala De Juárez"
$ws.Range("B50").Value = "Pachuca De Soto"
$ws.Range("B56").Value = "Autlán De Navarro"
$ws.Range("B88").Value = "Landa De Matamoros"
$ws.Range("B105").Value = "Paso De Ovejas"

# --- Drop trailing footnote/metadata rows (114-118), shrinking dimension to A1:D112 ---
$ws.Rows("114:118").Delete()
